# Add Excel "get list" functionality: replace the CSC440 grade-roster sheet
# (StudentID/CoursePrefix/CourseNum/Grade/Year/Semester) with a simpler
# Name/ID/Grade roster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- wipe the old 6-column table -----------------------------------------
$ws.Range("A1:F6").ClearContents()

# --- headers ---------------------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "ID"
$ws.Range("C1").Value = "Grade"

# --- data rows ---------------------------------------------------------------
$ws.Range("A2").Value = "Alex Hunter"
$ws.Range("B2").Value = 1111
$ws.Range("C2").Value = "A"

$ws.Range("A3").Value = "Jacob Anderson"
$ws.Range("B3").Value = 2222
$ws.Range("C3").Value = "C"

$ws.Range("A4").Value = "Mary Handerson"
$ws.Range("B4").Value = 3333
$ws.Range("C4").Value = "B"

$ws.Range("A5").Value = "Georeg Alan"
$ws.Range("B5").Value = 4444
$ws.Range("C5").Value = "F"

# --- center-align the ID and Grade columns (header + data) -----------------
$ws.Range("B1:C5").HorizontalAlignment = -4108

# --- column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 22.42578125
$ws.Columns.Item(2).ColumnWidth = 12.5703125
$ws.Columns.Item(3).ColumnWidth = 6.28515625

# --- match the saved selection/cursor position ------------------------------
$ws.Range("E3").Select() | Out-Null
